# "Find bugs summary" (#57) -- fix the Anesthesiology_Physician column header on the
# second sheet ("Dec 02 2020 - Dec 15 2020") to read "Anesthesiologist_Physician",
# and leave the sheet scrolled/selected the way the author left it (I4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dec 02 2020 - Dec 15 2020")

# Column I header text: Anesthesiology_Physician -> Anesthesiologist_Physician
$ws.Range("I1").Value = "Anesthesiologist_Physician"

# Make this the active sheet/window and restore the author's view state.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I4").Select()
